# Auto-generated Excel COM-interop script to update cryptos list
# Updates price (D) and volume-change (E) values, and a few coin name/link swaps
# reflecting the GitHub Actions refreshed cryptos.xlsx snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "66.977.91"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "2.479.36"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").Value = "2.478.99"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "2.923.41"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").Value = "66.789.47"
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Value = "2.464.92"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  -5.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.43%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D30").Value = "0.0₃0932"
$ws.Range("E30").Value = "  -6.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "510.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("E34").Value = "  -5.32%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.334"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.538"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.60%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0267"
$ws.Range("E50").Value = "  -4.14%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.64%  "
